# Update the "Soal" (Questions) sheet:
#  - Row 3 (Pilgan / multiple-choice single): fix the weight column E3 from
#    "5,0,0,0" to "25,0,0,0" and label its usage in column G.
#  - Row 4 (Essay): label its usage in column G.
#  - Row 5 (Praktikum): bump the weight in F5 from 10 to 25 and label its
#    usage in column G.
#  - Row 6 (new): add a second "Pilgan" (multiple-choice, multi-answer /
#    "Jamak") example question with its own answer choices and weights,
#    labelled in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Digunakan untuk Soal Pilihan Ganda Tipe Tunggal"

$ws.Range("B6").Value = "Manakah dari berikut ini yang TIDAK termasuk jenis algoritma Machine Learning?"
$ws.Range("C6").Value = "Pilgan"
$ws.Range("G6").Value = "Digunakan untuk Soal Pilihan Ganda Tipe Jamak"

$ws.Range("G5").Value = "Digunakan untuk Soal Praktikum"
$ws.Range("F5").Value = 25

$ws.Range("G4").Value = "Digunakan untuk Soal Essai"

$ws.Range("D6").Value = "Basis Data,PRG  5,PRG 6, Psikologi"
$ws.Range("E6").Value = "25,0,0,10"

$ws.Range("E3").Value = "25,0,0,0"
